$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, shifting rows 234:323 down to 235:324
$ws.Rows("234:234").Insert()

# Populate the new row 234 with its data
$ws.Range("A234").Value = 9
$ws.Range("B234").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C234").Value = "Metropolitana"
$ws.Range("D234").Value = 45006
$ws.Range("D234").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E234").Value = 13
$ws.Range("F234").Value = "Fruta"
$ws.Range("G234").Value = 100101
$ws.Range("H234").Value = "Berries"
$ws.Range("I234").Value = 100101001
$ws.Range("J234").Value = "Arándano (blue)"
$ws.Range("K234").Value = "Sin especificar"
$ws.Range("L234").Value = "Primera"
$ws.Range("M234").Value = 570
$ws.Range("N234").Value = 3800
$ws.Range("O234").Value = 4000
$ws.Range("P234").Value = 3898
$ws.Range("Q234").Value = "`$/bandeja 2 kilos"
$ws.Range("R234").Value = "Provincia de Curicó"
$ws.Range("S234").Value = 1949
$ws.Range("T234").Value = 2
